$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Buy CD - In Stock")
$ps = $ws1.PageSetup
$members = Get-Member -InputObject $ps
($members | ForEach-Object { $_.Name }) -join "`n"
